$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "shardul"
$ws.Range("B4").Value = "shardul123"
$ws.Range("A6").Value = "kiran"
$ws.Range("B6").Value = "kiran12"

$ws.Range("B13").Select()
